$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    4  = -11
    5  = -4
    13 = -5
    15 = -4
    17 = -5
    20 = -3
    23 = -2
    28 = -7
    29 = -2
    35 = -6
    36 = 5
    41 = -2
    46 = -1
    47 = -4
    48 = 1
    49 = -2
    50 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
